# edit.ps1 - apply spell-check (proofErr) run-splitting and add new
# "Hyperledger/IPFS/..." list items, per the supplied diff.

function Set-RunsXml {
    # Replaces the contents of $rng with the supplied run-level OOXML
    # ($innerXml is everything that belongs *inside* a <w:p> element -
    # runs, proofErr markers, etc). The surrounding <w:p> is NOT part of
    # $innerXml so that InsertXML only swaps the run content and leaves
    # the existing paragraph mark / pPr untouched.
    param($rng, $innerXml)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

function Replace-TextWithRuns {
    # Finds the first occurrence of $searchText anywhere in the document
    # and replaces just that span with $innerXml (run-level OOXML).
    param($d, $searchText, $innerXml)
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $searchText
        return
    }
    # IMPORTANT: build a brand-new Range from the Find-match bounds
    # instead of calling InsertXML directly on the (Find-mutated) $rng -
    # doing so avoids a leftover-text artifact left behind by Find.
    $clean = $d.Range($rng.Start, $rng.End)
    Set-RunsXml $clean $innerXml
}

$d = $word.ActiveDocument
$lang = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

# ---------------------------------------------------------------------
# 1) "Blockchain devs " -> "Blockchain " + proofErr(devs) + " "
# ---------------------------------------------------------------------
$inner1 = '<w:r>' + $lang + '<w:t xml:space="preserve">Blockchain </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>devs</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> </w:t></w:r>'
Replace-TextWithRuns $d "Blockchain devs " $inner1

# ---------------------------------------------------------------------
# 2) "...consensus protocols ,dApps, smart contracts..." -> split around dApps
# ---------------------------------------------------------------------
$inner2 = '<w:r>' + $lang + '<w:t>Those who optimizes and building blockchain related application like consensus protocols ,</w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>dApps</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t>, smart contracts ,design architecture systems.</w:t></w:r>'
Replace-TextWithRuns $d "Those who optimizes and building blockchain related application like consensus protocols ,dApps, smart contracts ,design architecture systems." $inner2

# ---------------------------------------------------------------------
# 3) " , Datastructure and Cryptography" -> split around Datastructure
# ---------------------------------------------------------------------
$inner3 = '<w:r>' + $lang + '<w:t xml:space="preserve"> , </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>Datastructure</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> and Cryptography</w:t></w:r>'
Replace-TextWithRuns $d " , Datastructure and Cryptography" $inner3

# ---------------------------------------------------------------------
# 3b) Insert 7 new list paragraphs right after the "...Cryptography" item
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("HTML , React , JavaScript", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $d.Range($rng.Start, $rng.Start).Paragraphs.Item(1)
$anchorIndex = $anchorPara.Range.Information(1)   # not reliable across runtimes - fallback below

# Robustly determine the paragraph index of the "...Cryptography" item by
# scanning the Paragraphs collection for the one ending in "Cryptography".
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*HTML , React , JavaScript*Cryptography*") {
        $idx = $i
        break
    }
}

$items = @("What is Hyperledger?", "What is IPFS?", "Cryptography?", "What is truffle ?", "What are NFT's ?", "What is DeFi?", "What is GETH ?")
foreach ($item in $items) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $d.Paragraphs($idx).Range.Text = $item
}

# Fix up "What is DeFi?" paragraph with a proofErr split around DeFi
$innerDefi = '<w:r>' + $lang + '<w:t xml:space="preserve">What is </w:t></w:r>' + `
             '<w:proofErr w:type="spellStart"/>' + `
             '<w:r>' + $lang + '<w:t>DeFi</w:t></w:r>' + `
             '<w:proofErr w:type="spellEnd"/>' + `
             '<w:r>' + $lang + '<w:t>?</w:t></w:r>'
Replace-TextWithRuns $d "What is DeFi?" $innerDefi

# ---------------------------------------------------------------------
# 4) hash function typo sentence -> split around inout / nay / aoutput / lengthof
# ---------------------------------------------------------------------
$inner4 = '<w:r>' + $lang + '<w:t xml:space="preserve">A hash function takes an </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>inout</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> of </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>nay</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> length and converts into an </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>aoutput</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> of a fixed </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>lengthof</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> 32 bytes.</w:t></w:r>'
Replace-TextWithRuns $d "A hash function takes an inout of nay length and converts into an aoutput of a fixed lengthof 32 bytes." $inner4

# ---------------------------------------------------------------------
# 5) "How to Miners and Transations in block?" -> split around Transations
# ---------------------------------------------------------------------
$inner5 = '<w:r>' + $lang + '<w:t xml:space="preserve">How to Miners and </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r>' + $lang + '<w:t>Transations</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r>' + $lang + '<w:t xml:space="preserve"> in block?</w:t></w:r>'
Replace-TextWithRuns $d "How to Miners and Transations in block?" $inner5

Write-Host "All edits applied."
